$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 86, shifting existing rows 86-130 down to 87-131.
$ws.Rows(86).Insert()

# Populate the newly inserted row with the new weekly observation.
$ws.Range("A86").Value = 7
$ws.Range("B86").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C86").Value = "Ñuble"
$ws.Range("D86").Value = 44466
$ws.Range("D86").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E86").Value = 16
$ws.Range("F86").Value = 100112032
$ws.Range("G86").Value = "Zapallo italiano"
$ws.Range("H86").Value = "Sin especificar"
$ws.Range("I86").Value = "Primera"
$ws.Range("J86").Value = 120
$ws.Range("K86").Value = 13000
$ws.Range("L86").Value = 14000
$ws.Range("M86").Value = 13500
$ws.Range("N86").Value = "$/caja 50 unidades"
$ws.Range("O86").Value = "Región de Arica y Parinacota"
$ws.Range("P86").Value = 270
$ws.Range("Q86").Value = 50
$ws.Range("R86").Value = "Hortaliza"
